# Generate Report for Handback
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: add a "zh-cn" results column and a data row for overview.md
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:B2"))
$loOverview.ShowHeaders = $true
$loOverview.ShowAutoFilter = $true

$wsOverview.Range("B1").Value = "zh-cn"
$wsOverview.Range("A2").Value = "overview.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "overview.md", "", "", "overview.md")
$wsOverview.Range("A2").Style = "HyperLink"
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# zh-cn sheet: add the handback-status detail row for overview.md
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:L2"))
$loZhCn.ShowHeaders = $true
$loZhCn.ShowAutoFilter = $true

$wsZhCn.Range("A2").Value = "overview.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "overview.md", "", "", "overview.md")
$wsZhCn.Range("A2").Style = "HyperLink"

$wsZhCn.Range("B2").Value = ".md"
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"

$wsZhCn.Range("D2").Value = "overview.8cc85dd99239e10c76baa6006d906abfd6122c3f.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), "overview.8cc85dd99239e10c76baa6006d906abfd6122c3f.zh-cn.xlf", "", "", "overview.8cc85dd99239e10c76baa6006d906abfd6122c3f.zh-cn.xlf")
$wsZhCn.Range("D2").Style = "HyperLink"

$wsZhCn.Range("E2").Value = "2016-04-12 05:41:03"
$wsZhCn.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Range("F2").Value = "overview.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "overview.md", "", "", "overview.md")
$wsZhCn.Range("F2").Style = "HyperLink"

$wsZhCn.Range("G2").Value = "overview.8cc85dd99239e10c76baa6006d906abfd6122c3f.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), "overview.8cc85dd99239e10c76baa6006d906abfd6122c3f.zh-cn.xlf", "", "", "overview.8cc85dd99239e10c76baa6006d906abfd6122c3f.zh-cn.xlf")
$wsZhCn.Range("G2").Style = "HyperLink"

$wsZhCn.Range("H2").Value = "2016-04-12 06:40:48"
$wsZhCn.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Range("J2").Value = "Include"
